# Fix Excel data-validation bug: the dropdown on user_department!C2:C100 was
# backed by a literal comma-separated string. Add a hidden lookup sheet
# ("departments") that holds the canonical department list and (per the
# upstream change) wire a list-type validation to it.

$wb  = $excel.ActiveWorkbook
$userDeptSheet = $wb.Worksheets.Item(1)

# --- create the "departments" lookup sheet -------------------------------
$deptSheet = $wb.Worksheets.Add()
$deptSheet.Name = "departments"

$departments = @(
    "Department",
    "Department of Commerce",
    "Department of Business Studies",
    "Department of Cultural Studies",
    "Department of Media Communications",
    "Department of Arts",
    "Department of Mathematics",
    "Department of Physical Education",
    "Department of Computer Science"
)

for ($i = 0; $i -lt $departments.Length; $i++) {
    $deptSheet.Range("A" + ($i + 1)).Value = $departments[$i]
}

$deptSheet.Columns.Item(1).ColumnWidth = 35.75

# list-validation on the new sheet's own C2:C100, sourced from the
# department list on this sheet (mirrors the upstream diff exactly)
$deptSheet.Range("C2:C100").Validation.Add(3, 1, 1, '=''departments''!$A$2:$A$9')

# hide the helper sheet and park it after "user_department" in tab order
$deptSheet.Visible = $false
$deptSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# re-resolve by name: the old reference goes stale across the Move() call
$deptSheet = $wb.Worksheets.Item("departments")

# the hidden "departments" sheet becomes the selected tab (matches diff:
# tabSelected moved off user_department's sheetView)
$deptSheet.Select()
